$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(137).Insert()

$ws.Range("A137").Value = 11
$ws.Range("B137").Value = "Vega Monumental Concepción"
$ws.Range("C137").Value = "Bíobío"
$ws.Range("D137").Value = 45009
$ws.Range("E137").Value = 8
$ws.Range("F137").Value = 100112043
$ws.Range("G137").Value = "Pepino ensalada"
$ws.Range("H137").Value = "Sin especificar"
$ws.Range("I137").Value = "Primera"
$ws.Range("J137").Value = 100
$ws.Range("K137").Value = 7000
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = 7250
$ws.Range("N137").Value = "`$/caja 60 unidades"
$ws.Range("O137").Value = "Región de Arica y Parinacota"
$ws.Range("P137").Value = 121
$ws.Range("Q137").Value = 60
$ws.Range("R137").Value = "Hortaliza"

Write-Host "done"
